# Sum_of_printed_Material.xlsx update
# - Inserts a new "Purpose" column (D) between Printer Used and Weight used(g)
# - Fills Purpose for every existing row
# - Updates row 6 (Bamboo/20*20 extrusion test) weight 20.6 -> 60.7
# - Updates row 8: material PLA (Yellow) -> PLA (Blue), printer -> Creality,
#   weight 273.31 -> 1000
# - Updates row 11: material PLA(Orange) -> PLA(blue), weight 342.5 -> 324
# - Resizes columns, sets column D to a bestFit-ish width
# - Moves selection to E19 (matches the saved cursor position in the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Weight used(g)" column (D),
# shifting it (and the trailing blank style column) one to the right.
$ws.Columns("D:D").Insert()

# Header for the new column
$ws.Range("D1").Value = "Purpose"

# --- Purpose values, filled in the same row-by-row order the sheet was
#     originally authored in (this controls the shared-string table order) ---
$ws.Range("D2").Value  = "Tolerance test of Extrusion"
$ws.Range("D3").Value  = "Tolerance test of Extrusion"
$ws.Range("D4").Value  = "Tolerance test of Extrusion"
$ws.Range("D5").Value  = "Pressfit Test"
$ws.Range("D6").Value  = "20*20 extrusion test"
$ws.Range("D7").Value  = "Servo motor case"
$ws.Range("D8").Value  = "Motor case print wasted"

# --- Row 8: airless-tyre motor-case entry gets corrected material/printer/weight ---
$ws.Range("B8").Value = "PLA (Blue)"
$ws.Range("C8").Value = "Creality"
$ws.Range("E8").Value = 1000

$ws.Range("D9").Value  = "Motor case printre used"
$ws.Range("D10").Value = "front Sus joint + servo Link"

# --- Row 11: material correction, then its Purpose ---
$ws.Range("B11").Value = "PLA(blue)"
$ws.Range("D11").Value = "threaded screws and joints"
$ws.Range("E11").Value = 324

$ws.Range("D12").Value = "small Joints + Links"

# --- Row 6: weight correction ---
$ws.Range("E6").Value = 60.7

# --- Column width: give the new Purpose column (D) a best-fit-style width.
#     A:C and E keep the 20.6328125 width inherited from the insert, and F
#     (the narrow trailing column) is untouched. ---
$ws.Columns("D:D").ColumnWidth = 22.42

# --- Restore the saved selection/cursor position ---
$ws.Range("E19").Select()
